$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each 12-row block (one calendar year) gets cyclically rotated so that
# Oct/Nov/Dec move to the top of the block and Jan..Sep follow.
$blockStarts = @(2, 14, 26, 38)

foreach ($base in $blockStarts) {
    # Snapshot all 12 rows x 10 cols (A..J) of this block before writing anything
    $snapshot = @()
    for ($i = 0; $i -lt 12; $i++) {
        $rowVals = @()
        for ($col = 1; $col -le 10; $col++) {
            $rowVals += $ws.Cells.Item($base + $i, $col).Value2
        }
        $snapshot += (,$rowVals)
    }

    # new offset i <- old offset (i + 9) % 12
    for ($i = 0; $i -lt 12; $i++) {
        $srcOffset = ($i + 9) % 12
        $srcRow = $snapshot[$srcOffset]
        for ($col = 1; $col -le 10; $col++) {
            $ws.Cells.Item($base + $i, $col).Value = $srcRow[$col - 1]
        }
    }
}